# FRA0400.xlsx - BOM pedalier update
# Adds the pedal-assembly parts (rails, pedals, supports, spacers, ...) to
# the "Pedals" BOM sheet, replacing the generic placeholder rows 3-8 with
# real part data and appending rows 9-15 for the remaining parts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Rows 3-8: fill in the real part name / make-or-buy / description /
# quantity for the parts that already had placeholder rows. The "ID"
# column (G) already holds the correct FR_0400x codes and is untouched.
# ---------------------------------------------------------------------
$ws.Range("C3").Value = "Upper rail"
$ws.Range("D3").Value = "m"
$ws.Range("E3").Value = "Short description of the part"
$ws.Range("F3").Value = 2

$ws.Range("C4").Value = "Lower rail"
$ws.Range("D4").Value = "m"
$ws.Range("E4").Value = "Short description of the part"
$ws.Range("F4").Value = 2

$ws.Range("C5").Value = "Accelerator pedal"
$ws.Range("D5").Value = "m"
$ws.Range("E5").Value = "Short description of the part"
$ws.Range("F5").Value = 1

$ws.Range("C6").Value = "Brake pedal"
$ws.Range("D6").Value = "m"
$ws.Range("E6").Value = "Short description of the part"
$ws.Range("F6").Value = 1

$ws.Range("C7").Value = "Foot top support"
$ws.Range("D7").Value = "m"
$ws.Range("E7").Value = "Short description of the part"
$ws.Range("F7").Value = 2

$ws.Range("C8").Value = "Heel support"
$ws.Range("D8").Value = "m"
$ws.Range("E8").Value = "Short description of the part"
$ws.Range("F8").Value = 2

# ---------------------------------------------------------------------
# Rows 9-15: brand new parts. Copy the formatting of row 8 (border /
# fill / font) down into each new row before writing the values so the
# new rows look like the rest of the table.
# ---------------------------------------------------------------------
$ws.Range("C8:G8").Copy()
$ws.Range("C9:G15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C9").Value = "Side support brake pedal"
$ws.Range("D9").Value = "m"
$ws.Range("E9").Value = "Short description of the part"
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = "FR_04007"

$ws.Range("C10").Value = "Left side support accelerator pedal"
$ws.Range("D10").Value = "m"
$ws.Range("E10").Value = "Short description of the part"
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = "FR_04008"

$ws.Range("C11").Value = "Right side support brake pedal"
$ws.Range("D11").Value = "m"
$ws.Range("E11").Value = "Short description of the part"
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = "FR_04009"

$ws.Range("C12").Value = "Brake over-travel switch support"
$ws.Range("D12").Value = "m"
$ws.Range("E12").Value = "Short description of the part"
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = "FR_04010"

$ws.Range("C13").Value = "Cable support axis"
$ws.Range("D13").Value = "m"
$ws.Range("E13").Value = "Axis to pull the accelerator cable"
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = "FR_04011"

$ws.Range("C14").Value = "Inside spacer master cylinder"
$ws.Range("D14").Value = "m"
$ws.Range("E14").Value = "To support the upper part of the master cylinder"
$ws.Range("F14").Value = 2
$ws.Range("G14").Value = "FR_04012"

$ws.Range("C15").Value = "Outside spacer master cylinder"
$ws.Range("D15").Value = "m"
$ws.Range("E15").Value = "To support the upper part of the master cylinder"
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = "FR_04013"

# ---------------------------------------------------------------------
# Row heights: short single-line rows stay close to the default data
# row height, rows with longer descriptions that wrap get auto-fit to a
# larger height by Excel.
# ---------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 18
$ws.Rows.Item(9).RowHeight = 15
$ws.Rows.Item(10).RowHeight = 27.6
$ws.Rows.Item(11).RowHeight = 27.6
$ws.Rows.Item(12).RowHeight = 27.6
$ws.Rows.Item(13).RowHeight = 15
$ws.Rows.Item(14).RowHeight = 27.6
$ws.Rows.Item(15).RowHeight = 27.6

# ---------------------------------------------------------------------
# View state: zoom to 70%, scroll so column C is leftmost, and leave the
# selection on F16 (just below the new last row), matching the author's
# final on-screen state.
# ---------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 70
$ws.Range("F16").Select()
$excel.ActiveWindow.ScrollColumn = 3
